# Apply updated crypto price/volume data to the "cryptos" worksheet.
# Source values are kept as literal text (matching original inline-string
# cells) so that numeric formatting/precision (trailing zeros, % suffix,
# scientific-looking small decimals) is preserved exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Range, $Text) {
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") "243.62"
Set-TextValue $ws.Range("E2") "-0.19%"

# Row 3
Set-TextValue $ws.Range("D3") "29.99"
Set-TextValue $ws.Range("E3") "13.37%"

# Row 4
Set-TextValue $ws.Range("D4") "5.153"
Set-TextValue $ws.Range("E4") "0.04%"

# Row 5
Set-TextValue $ws.Range("D5") "0.05671"
Set-TextValue $ws.Range("E5") "1.18%"

# Row 6
Set-TextValue $ws.Range("E6") "1.09%"

# Row 7
Set-TextValue $ws.Range("D7") "0.8478"
Set-TextValue $ws.Range("E7") "3.47%"

# Row 8
Set-TextValue $ws.Range("D8") "0.8607"
Set-TextValue $ws.Range("E8") "3.90%"

# Row 9
Set-TextValue $ws.Range("D9") "0.1350"
Set-TextValue $ws.Range("E9") "1.23%"

# Row 10
Set-TextValue $ws.Range("D10") "0.06925"
Set-TextValue $ws.Range("E10") "0.03%"

# Row 11
Set-TextValue $ws.Range("D11") "0.02893"
Set-TextValue $ws.Range("E11") "-0.16%"

# Row 12
Set-TextValue $ws.Range("D12") "0.09380"
Set-TextValue $ws.Range("E12") "-0.07%"

# Row 13
Set-TextValue $ws.Range("D13") "0.001514"
Set-TextValue $ws.Range("E13") "-0.10%"

# Row 14
Set-TextValue $ws.Range("D14") "0.04172"
Set-TextValue $ws.Range("E14") "-9.71%"

# Row 15
Set-TextValue $ws.Range("D15") "0.0005998"
Set-TextValue $ws.Range("E15") "-94.01%"

# Row 16
Set-TextValue $ws.Range("D16") "0.006086"
Set-TextValue $ws.Range("E16") "-2.55%"

# Row 17
Set-TextValue $ws.Range("E17") "-4.07%"

# Row 18
Set-TextValue $ws.Range("D18") "3.034"
Set-TextValue $ws.Range("E18") "0.05%"

# Row 19
Set-TextValue $ws.Range("D19") "2.134"
Set-TextValue $ws.Range("E19") "-2.26%"

# Row 21
Set-TextValue $ws.Range("D21") "0.03342"
Set-TextValue $ws.Range("E21") "8.16%"

# Row 22
Set-TextValue $ws.Range("E22") "0.29%"

# Row 23
Set-TextValue $ws.Range("D23") "3.634"
Set-TextValue $ws.Range("E23") "-2.88%"

# Row 24
Set-TextValue $ws.Range("E24") "2.40%"

# Row 25
Set-TextValue $ws.Range("D25") "0.001213"
Set-TextValue $ws.Range("E25") "-0.91%"

# Row 26
Set-TextValue $ws.Range("D26") "0.004441"
Set-TextValue $ws.Range("E26") "-1.09%"

# Row 27
Set-TextValue $ws.Range("D27") "0.0001180"

# Row 28
Set-TextValue $ws.Range("D28") "0.0001391"
Set-TextValue $ws.Range("E28") "-0.58%"

# Row 40
Set-TextValue $ws.Range("D40") "0.03738"
Set-TextValue $ws.Range("E40") "2.62%"

# Row 41
Set-TextValue $ws.Range("D41") "0.005848"
Set-TextValue $ws.Range("E41") "-5.18%"

# Row 42
Set-TextValue $ws.Range("D42") "0.1059"
Set-TextValue $ws.Range("E42") "0.82%"

# Row 43
Set-TextValue $ws.Range("D43") "0.002279"
Set-TextValue $ws.Range("E43") "-4.99%"

# Row 44
Set-TextValue $ws.Range("D44") "0.009286"
Set-TextValue $ws.Range("E44") "2.88%"

# Row 45
Set-TextValue $ws.Range("D45") "0.00005105"
Set-TextValue $ws.Range("E45") "-4.51%"

# Row 46
Set-TextValue $ws.Range("D46") "0.00000000750"
Set-TextValue $ws.Range("E46") "0.00%"

# Row 47
Set-TextValue $ws.Range("D47") "0.09995"
Set-TextValue $ws.Range("E47") "-30.56%"

# Row 48
Set-TextValue $ws.Range("D48") "0.002773"
Set-TextValue $ws.Range("E48") "-5.20%"

# Row 49
Set-TextValue $ws.Range("D49") "0.00002099"
Set-TextValue $ws.Range("E49") "0.00%"

# Row 50
Set-TextValue $ws.Range("D50") "0.0001999"
Set-TextValue $ws.Range("E50") "0.00%"
